# chore: update Sheets via scheduled runner
# Refreshes market-price-derived columns (currentAveragePrice / NQ / HQ,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) on a handful of leve rows across all
# eight class sheets, reflecting newly pulled marketboard data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 47.57143
$ws.Range("I9").Value = 50.666668
$ws.Range("J9").Value = 29
$ws.Range("K9").Value = 50.666668
$ws.Range("L9").Value = 29
$ws.Range("M9").Value = 118.333332
$ws.Range("N9").Value = -367

# Row 34
$ws.Range("H34").Value = 8995.5
$ws.Range("I34").Value = 8995.5
$ws.Range("K34").Value = 8995.5
$ws.Range("M34").Value = -8792.5

# Row 36
$ws.Range("H36").Value = 8995.5
$ws.Range("I36").Value = 8995.5
$ws.Range("K36").Value = 8995.5
$ws.Range("M36").Value = -8280.5

# Row 38
$ws.Range("H38").Value = 3747.1428
$ws.Range("I38").Value = 96.666664
$ws.Range("K38").Value = 289.999992
$ws.Range("M38").Value = 82.00000799999998

# Row 76
$ws.Range("H76").Value = 916.6667
$ws.Range("I76").Value = 916.6667
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 916.6667
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -601.6667
$ws.Range("N76").ClearContents()

# Row 79
$ws.Range("H79").Value = 916.6667
$ws.Range("I79").Value = 916.6667
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 916.6667
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = 175.3333
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 582.5185
$ws.Range("I32").Value = 590.3077
$ws.Range("K32").Value = 590.3077
$ws.Range("M32").Value = -303.3077

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2814.3333
$ws.Range("I20").Value = 2471.75
$ws.Range("J20").Value = 3499.5
$ws.Range("K20").Value = 2471.75
$ws.Range("L20").Value = 3499.5
$ws.Range("M20").Value = -2224.75
$ws.Range("N20").Value = -3993.5

# Row 105
$ws.Range("H105").Value = 1618.2858
$ws.Range("I105").Value = 1579.5
$ws.Range("K105").Value = 1579.5
$ws.Range("M105").Value = 167.5

$ws = $wb.Worksheets.Item("CRP")
# Row 2
$ws.Range("H2").Value = 246.5
$ws.Range("I2").Value = 188
$ws.Range("J2").Value = 305
$ws.Range("K2").Value = 188
$ws.Range("L2").Value = 305
$ws.Range("M2").Value = -75
$ws.Range("N2").Value = -531

# Row 22
$ws.Range("H22").Value = 2499
$ws.Range("I22").Value = 2415.8333
$ws.Range("K22").Value = 2415.8333
$ws.Range("M22").Value = -2065.8333

# Row 122
$ws.Range("H122").Value = 1255.375
$ws.Range("I122").Value = 1034.6666
$ws.Range("K122").Value = 3103.9998
$ws.Range("M122").Value = -653.9998000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 180.2
$ws.Range("I7").Value = 67
$ws.Range("J7").Value = 350
$ws.Range("K7").Value = 201
$ws.Range("L7").Value = 1050
$ws.Range("M7").Value = -89
$ws.Range("N7").Value = -1274

# Row 55
$ws.Range("H55").Value = 6883.25
$ws.Range("J55").Value = 6883.25
$ws.Range("L55").Value = 20649.75
$ws.Range("N55").Value = -21003.75

# Row 107
$ws.Range("H107").Value = 596.6
$ws.Range("I107").Value = 167
$ws.Range("J107").Value = 883
$ws.Range("K107").Value = 501
$ws.Range("L107").Value = 2649
$ws.Range("M107").Value = 1419
$ws.Range("N107").Value = -6489

# Row 117
$ws.Range("H117").Value = 537.3
$ws.Range("I117").Value = 69.666664
$ws.Range("K117").Value = 208.999992
$ws.Range("M117").Value = 3233.000008

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

# Row 70
$ws.Range("H70").Value = 2000
$ws.Range("I70").Value = 2000
$ws.Range("K70").Value = 2000
$ws.Range("M70").Value = -1730

# Row 73
$ws.Range("H73").Value = 2000
$ws.Range("I73").Value = 2000
$ws.Range("K73").Value = 2000
$ws.Range("M73").Value = -1064

# Row 80
$ws.Range("H80").Value = 2872.2222
$ws.Range("I80").Value = 3051.6667
$ws.Range("J80").Value = 2782.5
$ws.Range("K80").Value = 3051.6667
$ws.Range("L80").Value = 2782.5
$ws.Range("M80").Value = -2053.6667
$ws.Range("N80").Value = -4778.5

# Row 83
$ws.Range("H83").Value = 2872.2222
$ws.Range("I83").Value = 3051.6667
$ws.Range("J83").Value = 2782.5
$ws.Range("K83").Value = 15258.3335
$ws.Range("L83").Value = 13912.5
$ws.Range("M83").Value = -10266.3335
$ws.Range("N83").Value = -23896.5

# Row 113
$ws.Range("H113").Value = 5489.385
$ws.Range("I113").Value = 4467.875
$ws.Range("K113").Value = 4467.875
$ws.Range("M113").Value = -2297.875

# Row 132
$ws.Range("H132").Value = 2003.375
$ws.Range("I132").Value = 2002.4286
$ws.Range("K132").Value = 6007.2858
$ws.Range("M132").Value = -3477.2858

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 750
$ws.Range("I46").Value = 750
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 750
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -562
$ws.Range("N46").ClearContents()

# Row 68
$ws.Range("H68").Value = 8125
$ws.Range("I68").Value = 2500
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 2500
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -1751
$ws.Range("N68").Value = -11498

# Row 71
$ws.Range("H71").Value = 8125
$ws.Range("I71").Value = 2500
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 12500
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -8756
$ws.Range("N71").Value = -57488

# Row 132
$ws.Range("H132").Value = 3702.6924
$ws.Range("I132").Value = 3702.6924
$ws.Range("K132").Value = 11108.0772
$ws.Range("M132").Value = -8578.0772

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 10700
$ws.Range("I2").Value = 12875
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 12875
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -12763
$ws.Range("N2").Value = -2224

# Row 4
$ws.Range("H4").Value = 9340.817999999999
$ws.Range("I4").Value = 14421.429
$ws.Range("J4").Value = 449.75
$ws.Range("K4").Value = 14421.429
$ws.Range("L4").Value = 449.75
$ws.Range("M4").Value = -14308.429
$ws.Range("N4").Value = -675.75

# Row 7
$ws.Range("H7").Value = 6167
$ws.Range("I7").Value = 500
$ws.Range("K7").Value = 500
$ws.Range("M7").Value = -387

# Row 107
$ws.Range("H107").Value = 354.7647
$ws.Range("I107").Value = 354.7647
$ws.Range("K107").Value = 1064.2941
$ws.Range("M107").Value = 855.7058999999999
